$d = $word.ActiveDocument

# Locate the "MM/DD/YYYY" placeholder text (highlighted yellow) inside the
# "... underwent a Heart MRI examination on MM/DD/YYYY. Upon examination, I report ..."
# paragraph, so we can replace it with a hidden merge-field paragraph
# ("`r cmr_date_time`"), matching the pattern already used elsewhere in the
# template (e.g. "`r pronoun`", "`r extracardiac_incidental_describe`").

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("MM/DD/YYYY", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'MM/DD/YYYY' placeholder text"
}

$startPos = $findRange.Start
$endPos = $findRange.End

# Remove the "MM/DD/YYYY" runs entirely.
$delRange = $d.Range($startPos, $endPos)
$delRange.Text = ""

# Split the paragraph at that point: everything up to here stays in the
# original paragraph, everything after (the trailing ". Upon examination, I
# report ") moves into a new paragraph that inherits the same paragraph
# formatting.
$splitRange = $d.Range($startPos, $startPos)
[void]$splitRange.InsertParagraphAfter()

# Insert a brand-new hidden paragraph between the two halves holding the
# merge field "`r cmr_date_time`" (with the proofing-error wrappers Word
# adds around the field name, just like the surrounding fields).
$fieldXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:vanish/><w:specVanish/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">`r </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cmr_date_time</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>`</w:t></w:r></w:p>'

$insertRange = $d.Range($startPos, $startPos)
[void]$insertRange.InsertXML($fieldXml)
